$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - B2 changes from "Coding Ninja" to "CN"
$ws.Range("B2").Value = "CN"

# Row 3 - B3 and C3 change
$ws.Range("B3").Value = "GFG/CN"
$ws.Range("C3").Value = "Stickler Thief/Maximum sum of non-adjacent elements/House Robber I"

# Row 4 - new data
$ws.Range("A4").Value = 215
$ws.Range("B4").Value = "LC/CN"
$ws.Range("C4").Value = "House Robber II-Circular placement"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "Java"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "DP(Space optimizatio)"

# Update selection to E4
$ws.Range("E4").Select()
